$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title (appears twice: heading and bold run near the end)
Replace-Text "Play Lancelot slot game for free - review" "Play Lancelot for Free - A Slot Game Inspired by King Arthur"

# "What we like" bullet list
Replace-Text "High-quality graphics and attention to detail." "Inspired by the legends of King Arthur and the Knights of the Round Table"
Replace-Text "Free spins can be activated with stacked Wild symbols." "Players can select the lines and bet amount"
Replace-Text "Autoplay mode allows for automatic plays." "Autoplay mode and quick bet formula for convenience"
Replace-Text "Quick bet formula makes it easy to set betting amounts." "Can activate up to 50 free spins with stacked Wild symbols"

# "What we don't like" bullet list
Replace-Text "No bonus game available." "No bonus game"
Replace-Text "Sound effects can be a bit lackluster." "Lackluster sound"

# Closing italic summary paragraph
Replace-Text "Discover the world of King Arthur and the Knights of the Round Table with Lancelot online slot game. Play now for free and enjoy high-quality graphics." "Play Lancelot for free and experience the legends of King Arthur and the Knights of the Round Table in this online slot game."
